# edit.ps1
# Applies the "third option for final agreement" change:
# Adds a new scenario row for the Final WTO Fisheries Subsidies Agreement where the
# high-seas (HS) prohibition only applies to non-tuna vessels (i.e. excludes longline
# and purse seine / tuna gear types), positioned between the existing "complete HS
# prohibition" row and the "no HS prohibition" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Duplicate row 52 (the "complete HS prohibition" scenario) and
#    insert it as the new row 54, pushing the old rows 54 and 55 down
#    to rows 55 and 56 respectively.
# ------------------------------------------------------------------
$ws.Rows("52:52").Copy()
$ws.Rows("54:54").Insert()
$ws.Rows("54:54").RowHeight = 297

# ------------------------------------------------------------------
# 2. Row 51 - drop the stray AM51 value (not applicable to this row)
# ------------------------------------------------------------------
$ws.Range("AM51").Clear()

# ------------------------------------------------------------------
# 3. Row 52 - relabel as the "complete HS prohibition" variant and
#    drop the stray AM52 value
# ------------------------------------------------------------------
$ws.Range("E52").Value = "Final Fisheries Subsidies Agreement [complete HS prohibition]"
$ws.Range("G52").Value = "WT/MIN(22)/W/22 | Complete prohibition for HS fishing"
$ws.Range("P52").Value = "<ol><li>As written, this text would prohibit all capacity-enhancing and ambiguous subsidies as defined by Sumaila et al. (2019) to fishing in areas beyond the subsidizing Member's national jurisdiction and outside the competence of regional fisheries management organizations or agreements (RFMO/As).</li><li>Very few areas of the ocean are not under the jurisdiction of a RFMO/A, but it is difficult to determine on a global scale wehther vessels are fishing for species governed by those RFMO/As at any given point in time.</li><li>We therefore assume that all capacity-enhancing and ambiguous subsidies as defined by Sumaila et al. (2019) are prohibited to vessels fishing in areas beyond Members' national jurisdictions. By default, we consider vessels spending at least 5% of their total annual effort fishing on the high seas to be affected. We note that this may be an amibitious interpretation of this prohibition.</li><li>For the purposes of modeling, the rules on reflagging and unassessed stocks are not considered.</ol>"
$ws.Range("AM52").Clear()

# ------------------------------------------------------------------
# 4. Row 53 - unhide the "no HS prohibition" variant and drop the
#    stray AM53 value
# ------------------------------------------------------------------
$ws.Range("A53").Value = "Yes"
$ws.Range("AM53").Clear()

# ------------------------------------------------------------------
# 5. Row 54 (new) - turn the copy of row 52 into the new third option:
#    HS prohibition applies only to non-tuna vessels
# ------------------------------------------------------------------
$ws.Range("C54").Value = "Final WTO Fisheries Subsidies Agreement - HS portion only applies to vessels that are not longline or purse seine"
$ws.Range("E54").Value = "Final Fisheries Subsidies Agreement [non-tuna HS prohibition]"
$ws.Range("G54").Value = "WT/MIN(22)/W/22 | Non-tuna prohibition for HS fishing"
$ws.Range("P54").Value = "<ol><li>As written, this text would prohibit all capacity-enhancing and ambiguous subsidies as defined by Sumaila et al. (2019) to fishing in areas beyond the subsidizing Member's national jurisdiction and outside the competence of regional fisheries management organizations or agreements (RFMO/As).</li><li>Very few areas of the ocean are not under the jurisdiction of a RFMO/A, but it is difficult to determine on a global scale whether vessels are fishing for species governed by those RFMO/As at any given point in time.</li><li>We therefore assume that all capacity-enhancing an ambiguous subsidies as defined by Sumaila et al. (2019) are prohibited to non-tuna vessels fishing in areas beyond Members' national jurisdictions. We assume vessels utilizing the following gear types to be tuna fishing vessels: purse seines, drifting longlines, and pole and line. By default, we consider non-tuna fishing vessels spending at least 5% of their total annual effort fishing on the high seas to be affected. </li><li>For the purposes of modeling, the rules on reflagging and unassessed stocks are not considered.</ol>"
$ws.Range("AY54").Value = "HS/TUNA"
$ws.Range("AM54").Clear()

# ------------------------------------------------------------------
# 6. Update the view selection to match the author's saved state
# ------------------------------------------------------------------
$ws.Range("P54").Select()
